$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2:E51").NumberFormat = "@"

$ws.Cells.Item(2, 4).Value = '27.396.92'
$ws.Cells.Item(2, 5).Value = '  -0.63%  '
$ws.Cells.Item(3, 4).Value = '1.654.50'
$ws.Cells.Item(3, 5).Value = '  -0.30%  '
$ws.Cells.Item(4, 5).Value = '  -0.30%  '
$ws.Cells.Item(5, 4).Value = '213.43'
$ws.Cells.Item(5, 5).Value = '  -0.66%  '
$ws.Cells.Item(6, 5).Value = '  -0.31%  '
$ws.Cells.Item(7, 5).Value = '  -0.43%  '
$ws.Cells.Item(8, 4).Value = '23.51'
$ws.Cells.Item(8, 5).Value = '  +1.08%  '
$ws.Cells.Item(9, 5).Value = '  -0.51%  '
$ws.Cells.Item(10, 5).Value = '  -1.01%  '
$ws.Cells.Item(11, 4).Value = '0.0876'
$ws.Cells.Item(11, 5).Value = '  -0.37%  '
$ws.Cells.Item(12, 4).Value = '1.888.70'
$ws.Cells.Item(13, 4).Value = '1.651.53'
$ws.Cells.Item(13, 5).Value = '  -0.44%  '
$ws.Cells.Item(14, 5).Value = '  -1.23%  '
$ws.Cells.Item(15, 4).Value = '0.570'
$ws.Cells.Item(15, 5).Value = '  +3.71%  '
$ws.Cells.Item(16, 4).Value = '65.58'
$ws.Cells.Item(16, 5).Value = '  -0.49%  '
$ws.Cells.Item(17, 4).Value = '27.396.92'
$ws.Cells.Item(17, 5).Value = '  -0.52%  '
$ws.Cells.Item(18, 4).Value = '232.11'
$ws.Cells.Item(18, 5).Value = '  -6.25%  '
$ws.Cells.Item(19, 5).Value = '  -0.53%  '
$ws.Cells.Item(20, 4).Value = '7.44'
$ws.Cells.Item(20, 5).Value = '  -0.48%  '
$ws.Cells.Item(21, 5).Value = '  -0.35%  '
$ws.Cells.Item(22, 5).Value = '  -2.22%  '
$ws.Cells.Item(23, 4).Value = '9.44'
$ws.Cells.Item(23, 5).Value = '  +3.81%  '
$ws.Cells.Item(24, 5).Value = '  +0.02%  '
$ws.Cells.Item(25, 5).Value = '  +0.70%  '
$ws.Cells.Item(26, 4).Value = '7.11'
$ws.Cells.Item(26, 5).Value = '  -0.73%  '
$ws.Cells.Item(27, 4).Value = '15.87'
$ws.Cells.Item(27, 5).Value = '  -2.06%  '
$ws.Cells.Item(28, 4).Value = '0.999'
$ws.Cells.Item(28, 5).Value = '  -0.56%  '
$ws.Cells.Item(29, 5).Value = '  +0.07%  '
$ws.Cells.Item(30, 5).Value = '  -0.49%  '
$ws.Cells.Item(31, 5).Value = '  -4.05%  '
$ws.Cells.Item(32, 5).Value = '  -1.46%  '
$ws.Cells.Item(33, 2).Value = 'Maker'
$ws.Cells.Item(33, 3).Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Cells.Item(33, 4).Value = '1.427.71'
$ws.Cells.Item(33, 5).Value = '  -0.24%  '
$ws.Cells.Item(34, 2).Value = 'InternetComputer(DFINITY)'
$ws.Cells.Item(34, 3).Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Cells.Item(34, 4).Value = '3.13'
$ws.Cells.Item(34, 5).Value = '  -0.15%  '
$ws.Cells.Item(35, 5).Value = '  +1.20%  '
$ws.Cells.Item(36, 5).Value = '  -1.07%  '
$ws.Cells.Item(37, 4).Value = '0.908'
$ws.Cells.Item(37, 5).Value = '  -2.35%  '
$ws.Cells.Item(38, 4).Value = '0.572'
$ws.Cells.Item(38, 5).Value = '  -1.09%  '
$ws.Cells.Item(39, 5).Value = '  +0.25%  '
$ws.Cells.Item(40, 5).Value = '  +0.15%  '
$ws.Cells.Item(41, 5).Value = '  -0.45%  '
$ws.Cells.Item(42, 4).Value = '5.53'
$ws.Cells.Item(42, 5).Value = '  +2.18%  '
$ws.Cells.Item(43, 2).Value = 'Aave'
$ws.Cells.Item(43, 3).Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Cells.Item(43, 4).Value = '65.01'
$ws.Cells.Item(43, 5).Value = '  -5.66%  '
$ws.Cells.Item(44, 2).Value = 'TrustWalletToken'
$ws.Cells.Item(44, 3).Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Cells.Item(44, 4).Value = '0.794'
$ws.Cells.Item(44, 5).Value = '  +0.80%  '
$ws.Cells.Item(45, 2).Value = 'MXToken'
$ws.Cells.Item(45, 3).Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Cells.Item(45, 4).Value = '2.22'
$ws.Cells.Item(45, 5).Value = '  +0.25%  '
$ws.Cells.Item(46, 4).Value = '1.796.93'
$ws.Cells.Item(46, 5).Value = '  -0.01%  '
$ws.Cells.Item(47, 4).Value = '1.69'
$ws.Cells.Item(47, 5).Value = '  -1.05%  '
$ws.Cells.Item(48, 4).Value = '88.02'
$ws.Cells.Item(48, 5).Value = '  -0.92%  '
$ws.Cells.Item(49, 4).Value = '0.0₆0106'
$ws.Cells.Item(49, 5).Value = '  +0.51%  '
$ws.Cells.Item(50, 5).Value = '  +0.01%  '
$ws.Cells.Item(51, 4).Value = '7.71'
$ws.Cells.Item(51, 5).Value = '  -0.80%  '

